$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ref, $val)
    $c = $ws.Range($ref)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "26.055.42"
Set-TextCell "E2" "  -1.13%  "
Set-TextCell "D3" "1.643.79"
Set-TextCell "E3" "  -1.47%  "
Set-TextCell "D5" "217.54"
Set-TextCell "E5" "  -0.93%  "
Set-TextCell "D6" "0.5189"
Set-TextCell "E6" "  -3.03%  "
Set-TextCell "E7" "  -0.58%  "
Set-TextCell "D8" "0.2614"
Set-TextCell "E8" "  -1.99%  "
Set-TextCell "D9" "0.06278"
Set-TextCell "E9" "  -2.03%  "
Set-TextCell "D10" "20.42"
Set-TextCell "E10" "  -2.37%  "
Set-TextCell "D11" "0.07750"
Set-TextCell "E11" "  -1.30%  "
Set-TextCell "B12" "Polkadot"
Set-TextCell "C12" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D12" "4.469"
Set-TextCell "E12" "  -2.26%  "
Set-TextCell "B13" "WrappedEther"
Set-TextCell "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D13" "1.643.81"
Set-TextCell "E13" "  -1.42%  "
Set-TextCell "D14" "1.869.33"
Set-TextCell "E14" "  -1.42%  "
Set-TextCell "D15" "0.5573"
Set-TextCell "E15" "  +0.53%  "
Set-TextCell "D16" "0.0₅7988"
Set-TextCell "E16" "  -2.62%  "
Set-TextCell "D17" "64.68"
Set-TextCell "E17" "  -1.89%  "
Set-TextCell "D18" "26.052.28"
Set-TextCell "E18" "  -1.22%  "
Set-TextCell "E19" "  -0.60%  "
Set-TextCell "D20" "4.620"
Set-TextCell "E20" "  -1.61%  "
Set-TextCell "D21" "192.59"
Set-TextCell "E21" "  -0.65%  "
Set-TextCell "E22" "  -2.60%  "
Set-TextCell "D23" "5.946"
Set-TextCell "E23" "  -1.81%  "
Set-TextCell "E24" "  -0.64%  "
Set-TextCell "D25" "146.40"
Set-TextCell "E25" "  -0.10%  "
Set-TextCell "D26" "0.1199"
Set-TextCell "E26" "  -2.79%  "
Set-TextCell "D27" "7.158"
Set-TextCell "E27" "  -0.85%  "
Set-TextCell "E28" "  -1.87%  "
Set-TextCell "D29" "1.480"
Set-TextCell "E29" "  -1.40%  "
Set-TextCell "D30" "0.05636"
Set-TextCell "E30" "  -4.08%  "
Set-TextCell "E31" "  -1.88%  "
Set-TextCell "D32" "3.451"
Set-TextCell "E32" "  -5.33%  "
Set-TextCell "D33" "3.352"
Set-TextCell "E33" "  +1.86%  "
Set-TextCell "D34" "1.595"
Set-TextCell "E34" "  -0.94%  "
Set-TextCell "D35" "2.789"
Set-TextCell "E35" "  -1.64%  "
Set-TextCell "E36" "  -0.43%  "
Set-TextCell "D37" "0.9371"
Set-TextCell "E37" "  -3.64%  "
Set-TextCell "E38" "  -3.30%  "
Set-TextCell "D39" "5.950"
Set-TextCell "E39" "  +1.78%  "
Set-TextCell "E40" "  -1.88%  "
Set-TextCell "D41" "1.050.58"
Set-TextCell "E41" "  -1.44%  "
Set-TextCell "E42" "  -0.64%  "
Set-TextCell "D43" "0.8402"
Set-TextCell "E43" "  -3.20%  "
Set-TextCell "D44" "102.30"
Set-TextCell "E44" "  -2.56%  "
Set-TextCell "D45" "1.780.25"
Set-TextCell "E45" "  -1.48%  "
Set-TextCell "D46" "56.84"
Set-TextCell "E46" "  -1.90%  "
Set-TextCell "B47" "Frax"
Set-TextCell "C47" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell "D47" "1.009"
Set-TextCell "E47" "  -0.59%  "
Set-TextCell "B48" "Cronos"
Set-TextCell "C48" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D48" "0.05304"
Set-TextCell "E48" "  +2.62%  "
Set-TextCell "B49" "Mantle"
Set-TextCell "C49" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D49" "0.4324"
Set-TextCell "E49" "  -1.42%  "
Set-TextCell "B50" "BabyDogeCoin"
Set-TextCell "C50" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D50" "0.0₈102"
Set-TextCell "E50" "  -4.12%  "
Set-TextCell "D51" "7.912"
Set-TextCell "E51" "  -0.90%  "

Write-Host "Applied 100 cell updates"
